$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force plain numeric-looking price strings to remain stored as text
# (matches the source data, which encodes "Price" as a text column).
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D14", "D15", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D42", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "40.674.23"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "2.373.51"
$ws.Range("E3").Value = "  -4.16%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "310.86"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("D6").Value = "87.30"
$ws.Range("E6").Value = "  -6.34%  "
$ws.Range("D7").Value = "0.528"
$ws.Range("E7").Value = "  -4.35%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -4.55%  "
$ws.Range("D10").Value = "0.0837"
$ws.Range("E10").Value = "  -4.97%  "
$ws.Range("D11").Value = "30.73"
$ws.Range("E11").Value = "  -7.94%  "
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "2.740.98"
$ws.Range("E13").Value = "  -4.14%  "
$ws.Range("D14").Value = "6.57"
$ws.Range("E14").Value = "  -5.40%  "
$ws.Range("D15").Value = "15.02"
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("D16").Value = "2.341.58"
$ws.Range("E16").Value = "  -5.54%  "
$ws.Range("D17").Value = "0.763"
$ws.Range("E17").Value = "  -4.60%  "
$ws.Range("D18").Value = "40.606.01"
$ws.Range("E18").Value = "  -2.72%  "
$ws.Range("D19").Value = "0.0₃0912"
$ws.Range("E19").Value = "  -4.31%  "
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  -4.88%  "
$ws.Range("D21").Value = "68.85"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").Value = "10.91"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").Value = "232.84"
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("D24").Value = "2.65"
$ws.Range("E24").Value = "  -4.13%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "1.81"
$ws.Range("E26").Value = "  -7.38%  "
$ws.Range("D27").Value = "23.83"
$ws.Range("E27").Value = "  -5.35%  "
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").Value = "9.39"
$ws.Range("E29").Value = "  -3.91%  "
$ws.Range("D30").Value = "33.82"
$ws.Range("E30").Value = "  -8.66%  "
$ws.Range("D31").Value = "152.74"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.26"
$ws.Range("E33").Value = "  -4.70%  "
$ws.Range("D34").Value = "0.0732"
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("E35").Value = "  -4.95%  "
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("D37").Value = "16.05"
$ws.Range("E37").Value = "  -8.23%  "
$ws.Range("D38").Value = "2.78"
$ws.Range("E38").Value = "  -4.98%  "
$ws.Range("D39").Value = "0.0997"
$ws.Range("E39").Value = "  -4.51%  "
$ws.Range("D40").Value = "1.72"
$ws.Range("E40").Value = "  -8.55%  "
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("D42").Value = "2.37"
$ws.Range("E42").Value = "  -6.03%  "
$ws.Range("D43").Value = "1.962.05"
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("D44").Value = "0.0271"
$ws.Range("E44").Value = "  -4.92%  "
$ws.Range("D45").Value = "17.70"
$ws.Range("E45").Value = "  -7.49%  "
$ws.Range("D46").Value = "9.57"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "2.73"
$ws.Range("E47").Value = "  -8.39%  "
$ws.Range("D48").Value = "2.613.48"
$ws.Range("E48").Value = "  -3.88%  "
$ws.Range("D49").Value = "93.56"
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("D50").Value = "72.52"
$ws.Range("E50").Value = "  -5.95%  "
$ws.Range("D51").Value = "50.81"
$ws.Range("E51").Value = "  -3.10%  "
